$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 3 de Abril de 2020 a las 16:20"

# --- Update per-country statistics (value-only refreshes, no reordering) ---

# Row 4: Estados Unidos
$ws.Range("E4").Value = 228932
$ws.Range("G4").Value = 29
$ws.Range("H4").Value = 6099

# Row 7: Alemania
$ws.Range("B7").Value = 87244
$ws.Range("C7").Value = 2450
$ws.Range("E7").Value = 61531
$ws.Range("G7").Value = 31
$ws.Range("H7").Value = 1138

# Row 16: Austria
$ws.Range("B16").Value = 11412
$ws.Range("C16").Value = 283
$ws.Range("E16").Value = 9222

# Row 20: Brasil
$ws.Range("B20").Value = 8165
$ws.Range("C20").Value = 121
$ws.Range("E20").Value = 7706
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = 332

# Row 24: Noruega
$ws.Range("E24").Value = 5208
$ws.Range("G24").Value = 6
$ws.Range("H24").Value = 56

# Row 31: Rumania
$ws.Range("E31").Value = 2778
$ws.Range("G31").Value = 7
$ws.Range("H31").Value = 122

# Row 42: Finlandia
$ws.Range("E42").Value = 1295
$ws.Range("G42").Value = 1
$ws.Range("H42").Value = 20

# Row 51: Argentina
$ws.Range("E51").Value = 970
$ws.Range("G51").Value = 3
$ws.Range("H51").Value = 39

# Row 74: Letonia
$ws.Range("E74").Value = 492
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 0

# Row 75: Bulgaria
$ws.Range("B75").Value = 485
$ws.Range("C75").Value = 28
$ws.Range("E75").Value = 441
$ws.Range("G75").Value = 4
$ws.Range("H75").Value = 14

# --- Mauricio overtakes Islas Feroe and Estado de Palestina in the ranking ---
# Mauricio's updated figures now place it right after Nigeria (row 104),
# pushing Islas Feroe (previously row 105) and Estado de Palestina
# (previously row 106) down one position each. Montenegro (row 108) is
# unaffected and keeps its data.

$ws.Range("A105").Value = "Mauricio"
$ws.Range("B105").Value = 186
$ws.Range("C105").Value = 17
$ws.Range("D105").Value = 0
$ws.Range("E105").Value = 179
$ws.Range("F105").Value = 1
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 7

$ws.Range("A106").Value = "Islas Feroe"
$ws.Range("B106").Value = 179
$ws.Range("C106").Value = 2
$ws.Range("D106").Value = 91
$ws.Range("E106").Value = 88
$ws.Range("F106").Value = 1
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 0

$ws.Range("A107").Value = "Estado de Palestina"
$ws.Range("B107").Value = 171
$ws.Range("C107").Value = 10
$ws.Range("D107").Value = 18
$ws.Range("E107").Value = 152
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 1

$wb.Save()
